# correct calling of run_optimization
# - Column A (Active flag): a handful of rows flip from TRUE to FALSE.
# - Column G: every data row's numeric 10000 becomes the text string "10000".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose "Active" boolean (column A) flips from TRUE to FALSE.
$rowsToDeactivate = @(4, 5, 6, 8, 9, 12, 13, 14, 16, 17, 21, 22, 24, 25, 27, 28, 30)
foreach ($r in $rowsToDeactivate) {
    $ws.Cells.Item($r, 1).Value = $false
}

# Column G (rows 2-30): store "10000" as text instead of a number.
for ($r = 2; $r -le 30; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $cell.NumberFormat = "@"
    $cell.Value = "10000"
    $cell.Style = "Normal"
}
